$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename Vada Pav -> वडापाव and add a hyperlink to its image on A2 ---
$ws.Hyperlinks.Add($ws.Range("A2"), "images/vadapav.avif", "", "", "वडापाव")

# --- Fix Cold Drink price (row 8, column B): 30 -> 20 ---
$ws.Range("B8").Value = 20

# --- Append new menu items (rows 9-16) ---
$newItems = @(
    @("Bhaje",             20,  "images/vadapav.avif"),
    @("Cake 500g",         300, "images/tea.png"),
    @("Cake 1000g",        600, "images/tea.png"),
    @("Kurkure",           10,  "images/coffee.png"),
    @("Pani bottle Small", 10,  "images/idli.png"),
    @("Pani Bottle Mothi", 20,  "images/dosa.png"),
    @("Manchurian",        40,  "images/panipuri.png"),
    @("Misal Pav",         70,  "images/cold.png")
)

$row = 9
foreach ($item in $newItems) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]

    $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 3)).WrapText = $true

    if ($item[0] -eq "Pani bottle Small" -or $item[0] -eq "Pani Bottle Mothi") {
        $ws.Rows.Item($row).RowHeight = 23.85
    } else {
        $ws.Rows.Item($row).RowHeight = 12.8
    }

    $row = $row + 1
}

# --- Column sizing: new narrower column A, slightly adjust column C ---
$ws.Columns.Item(1).ColumnWidth = 16.5
$ws.Columns.Item(3).ColumnWidth = 26.54

# --- Selection moves to A2 (where the hyperlink now lives) ---
$ws.Range("A2").Select()
